$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) values are stored as text in the sheet (many contain
# thousands separators like "69.609.90" that are not valid numbers).
# Force text format before assigning so plain-numeric-looking values
# (e.g. "672.48") are not auto-converted to numeric cells by Excel.

$r = $ws.Range("D2")
$r.NumberFormat = "@"
$r.Value = "69.609.90"
$r = $ws.Range("D3")
$r.NumberFormat = "@"
$r.Value = "3.705.86"
$ws.Range("E3").Value = "  +0.91%  "
$ws.Range("E4").Value = "  -0.19%  "
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "672.48"
$ws.Range("E5").Value = "  -1.33%  "
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "161.86"
$ws.Range("E6").Value = "  +2.53%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +1.08%  "
$ws.Range("E9").Value = "  +0.95%  "
$ws.Range("E10").Value = "  +1.90%  "
$ws.Range("E11").Value = "  +1.95%  "
$ws.Range("E12").Value = "  +1.51%  "
$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = "32.85"
$ws.Range("E13").Value = "  +2.20%  "
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = "3.698.31"
$ws.Range("E14").Value = "  +0.59%  "
$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = "69.655.76"
$ws.Range("E15").Value = "  +0.58%  "
$ws.Range("E16").Value = "  +1.83%  "
$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = "16.32"
$ws.Range("E17").Value = "  +2.78%  "
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = "474.26"
$ws.Range("E19").Value = "  +1.12%  "
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "9.81"
$ws.Range("E20").Value = "  -1.36%  "
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = "0.654"
$ws.Range("E21").Value = "  +1.13%  "
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = "80.39"
$ws.Range("E22").Value = "  +0.54%  "
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = "3.854.49"
$ws.Range("E24").Value = "  +6.19%  "
$ws.Range("E25").Value = "  -0.03%  "
$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = "11.01"
$ws.Range("E26").Value = "  +1.21%  "
$ws.Range("E27").Value = "  +0.74%  "
$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = "2.69"
$ws.Range("E28").Value = "  +0.02%  "
$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = "1.73"
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("E30").Value = "  +1.63%  "
$ws.Range("B31").Value = "Kaspa"
$ws.Range("C31").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = "0.168"
$ws.Range("E31").Value = "  +6.47%  "
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = "6.59"
$ws.Range("E32").Value = "  +0.74%  "
$ws.Range("E33").Value = "  +0.62%  "
$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = "0.997"
$ws.Range("E34").Value = "  -0.21%  "
$r = $ws.Range("D35")
$r.NumberFormat = "@"
$r.Value = "3.696.36"
$ws.Range("E35").Value = "  +1.16%  "
$ws.Range("E36").Value = "  +4.62%  "
$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = "6.12"
$ws.Range("E37").Value = "  +1.70%  "
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("E39").Value = "  +1.70%  "
$ws.Range("E40").Value = "  -0.08%  "
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = "0.0915"
$ws.Range("E41").Value = "  +1.97%  "
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = "173.93"
$ws.Range("E42").Value = "  +4.04%  "
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = "0.942"
$ws.Range("E43").Value = "  +0.24%  "
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = "47.08"
$ws.Range("E44").Value = "  -1.05%  "
$ws.Range("E45").Value = "  +2.51%  "
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = "0.000280"
$ws.Range("E46").Value = "  +1.74%  "
$ws.Range("E47").Value = "  +2.92%  "
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = "27.76"
$ws.Range("E48").Value = "  +3.63%  "
$ws.Range("E49").Value = "  -0.25%  "
$ws.Range("E50").Value = "  +1.92%  "
$ws.Range("E51").Value = "  +1.30%  "
